# "Add files via upload / updated worker page"
#
# The worksheet held a handful of throwaway test rows (a/b, ss/sss, aa/aaaaaa)
# sandwiched between real user records. This update removes that test data by
# deleting the three rows, which shifts the two genuine rows below them
# (shani/waizman and Michael/Elisha) up into rows 4 and 5, and updates the
# workbook's default font from Arial to Calibri.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook-wide default font: Arial -> Calibri.
[void]($wb.Styles.Item("Normal").Font.Name = "Calibri")

# Drop the three placeholder/test rows (old rows 4-6: a/b, ss/sss, aa/aaaaaa).
# Rows 7-8 (shani/waizman, Michael/Elisha) shift up to become rows 4-5.
$ws.Rows("4:6").Delete()

# Leave the selection where the editor left it.
[void]$ws.Range("I13").Select()
